# Hybrid_test_v2.xlsx: "Change network line settings (#14)"
# The R + L + (C // G) branch has been changed to (R + L)//C//G branch.
#
# On the "NetworkLine" sheet:
#  - rows 10-12 are the transformer-style "self" branches (1-2, 3-4, 5-6).
#    Their G (pu) column (F) used to be "inf" and is now 0.
#  - rows 13-18 are the per-bus shunt branches (1-1 .. 6-6).
#    Their R (pu) and wL (pu) columns (C, D) used to be 0 and are now "inf".

$wb = $excel.ActiveWorkbook

$wsLine = $wb.Worksheets.Item("NetworkLine")

$wsLine.Range("F10").Value = 0
$wsLine.Range("F11").Value = 0
$wsLine.Range("F12").Value = 0

$wsLine.Range("C13:D18").Value = "inf"

# NetworkLine becomes the active sheet/tab, with D14 as the selected cell
# (previously the Apparatus sheet was active and NetworkLine had B19 selected).
$wsLine.Activate()
$wsLine.Range("D14").Select()
